$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a pure number written as text in the source data;
# force text format so Excel does not auto-convert them to numeric values.
$textCells = @("D5", "D10", "D19", "D25", "D27", "D30", "D35", "D37", "D38", "D39", "D40", "D45", "D46", "D49", "D50", "D51")
foreach ($tc in $textCells) {
    $ws.Range($tc).NumberFormat = "@"
}

$ws.Range("D2").Value = '26.707.60'
$ws.Range("E2").Value = '  -0.35%  '
$ws.Range("D3").Value = '1.638.81'
$ws.Range("E4").Value = '  -0.27%  '
$ws.Range("D5").Value = '217.82'
$ws.Range("E5").Value = '  +0.50%  '
$ws.Range("E6").Value = '  -0.88%  '
$ws.Range("E7").Value = '  -0.26%  '
$ws.Range("E8").Value = '  -0.60%  '
$ws.Range("E9").Value = '  -0.92%  '
$ws.Range("D10").Value = '19.08'
$ws.Range("E10").Value = '  -0.58%  '
$ws.Range("E11").Value = '  -0.07%  '
$ws.Range("E12").Value = '  -0.75%  '
$ws.Range("D13").Value = '1.617.57'
$ws.Range("E13").Value = '  -2.10%  '
$ws.Range("E14").Value = '  -1.34%  '
$ws.Range("E15").Value = '  -1.65%  '
$ws.Range("E16").Value = '  -1.63%  '
$ws.Range("D17").Value = '26.691.15'
$ws.Range("E17").Value = '  -0.45%  '
$ws.Range("E18").Value = '  -2.49%  '
$ws.Range("D19").Value = '211.13'
$ws.Range("E19").Value = '  -3.36%  '
$ws.Range("E20").Value = '  -0.26%  '
$ws.Range("E21").Value = '  -0.63%  '
$ws.Range("E22").Value = '  -1.36%  '
$ws.Range("E23").Value = '  -3.25%  '
$ws.Range("E24").Value = '  -2.80%  '
$ws.Range("D25").Value = '146.72'
$ws.Range("E25").Value = '  -0.04%  '
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").Value = '0.117'
$ws.Range("E27").Value = '  -2.06%  '
$ws.Range("E28").Value = '  -0.76%  '
$ws.Range("E29").Value = '  -1.10%  '
$ws.Range("D30").Value = '0.0502'
$ws.Range("E30").Value = '  -3.06%  '
$ws.Range("E31").Value = '  +0.25%  '
$ws.Range("E32").Value = '  +0.07%  '
$ws.Range("E33").Value = '  -0.88%  '
$ws.Range("D34").Value = '1.266.65'
$ws.Range("E34").Value = '  -1.59%  '
$ws.Range("D35").Value = '1.53'
$ws.Range("E35").Value = '  -1.18%  '
$ws.Range("E36").Value = '  -0.90%  '
$ws.Range("D37").Value = '0.0175'
$ws.Range("E37").Value = '  -2.29%  '
$ws.Range("D38").Value = '0.527'
$ws.Range("E38").Value = '  -2.03%  '
$ws.Range("B39").Value = 'ARBITRUM'
$ws.Range("C39").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D39").Value = '0.804'
$ws.Range("E39").Value = '  -3.03%  '
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").Value = '1.01'
$ws.Range("E40").Value = '  -0.33%  '
$ws.Range("E41").Value = '  -1.37%  '
$ws.Range("E42").Value = '  -3.37%  '
$ws.Range("E43").Value = '  -3.73%  '
$ws.Range("E44").Value = '  -0.79%  '
$ws.Range("D45").Value = '91.37'
$ws.Range("E45").Value = '  -0.69%  '
$ws.Range("D46").Value = '60.08'
$ws.Range("E46").Value = '  +0.61%  '
$ws.Range("E47").Value = '  -1.97%  '
$ws.Range("E48").Value = '  +0.34%  '
$ws.Range("D49").Value = '7.52'
$ws.Range("E49").Value = '  -3.01%  '
$ws.Range("D50").Value = '0.0960'
$ws.Range("E50").Value = '  -1.11%  '
$ws.Range("B51").Value = 'USDD'
$ws.Range("C51").Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range("D51").Value = '1.01'
$ws.Range("E51").Value = '  -0.25%  '
